# Updates the RKI COVID-19 death-count workbook with a newer data pull
# (weekly sheet "COVID_Todesfälle" + monthly sheet "COVID_Todesfälle_Monat").
#
# Values are stored in the sheet as TEXT (shared strings) even though they
# look numeric, so every numeric-looking value is entered with a leading
# apostrophe to force text entry, then the cell style is reset back to
# "Normal" so no stray number-format survives on the cell itself.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force a (possibly numeric-looking) string to be stored as text,
    # mirroring the existing shared-string cells in this workbook.
    $needsQuote = $text -match '^-?[0-9]+(\.[0-9]+)?$'
    if ($needsQuote) {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

# ---------------------------------------------------------------------
# Sheet 1: COVID_Todesfälle (weekly counts)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Set-TextValue $ws1.Cells.Item(9, 2)  "<4"
Set-TextValue $ws1.Cells.Item(10, 2) "<4"

Set-TextValue $ws1.Cells.Item(18, 2) "1598"
Set-TextValue $ws1.Cells.Item(19, 2) "1170"
Set-TextValue $ws1.Cells.Item(20, 2) "783"
Set-TextValue $ws1.Cells.Item(21, 2) "515"
Set-TextValue $ws1.Cells.Item(22, 2) "352"
Set-TextValue $ws1.Cells.Item(23, 2) "272"

Set-TextValue $ws1.Cells.Item(25, 2) "113"
Set-TextValue $ws1.Cells.Item(26, 2) "72"

Set-TextValue $ws1.Cells.Item(30, 2) "27"
Set-TextValue $ws1.Cells.Item(33, 2) "29"

Set-TextValue $ws1.Cells.Item(40, 2) "54"

Set-TextValue $ws1.Cells.Item(43, 2) "117"
Set-TextValue $ws1.Cells.Item(44, 2) "232"
Set-TextValue $ws1.Cells.Item(45, 2) "389"
Set-TextValue $ws1.Cells.Item(46, 2) "756"
Set-TextValue $ws1.Cells.Item(47, 2) "1158"
Set-TextValue $ws1.Cells.Item(48, 2) "1527"
Set-TextValue $ws1.Cells.Item(49, 2) "1940"

# New week 48 row, appended at the bottom.
$ws1.Cells.Item(50, 1).Value = 48
Set-TextValue $ws1.Cells.Item(50, 2) "2579"

# ---------------------------------------------------------------------
# Sheet 2: COVID_Todesfälle_Monat (monthly counts)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# A new month (Feb, "<4" deaths) appears before the existing data, so
# insert a row and push everything else down.
$ws2.Rows.Item(2).Insert()
$ws2.Cells.Item(2, 1).Value = 2
Set-TextValue $ws2.Cells.Item(2, 2) "<4"

# Refreshed counts for the existing months (now shifted down one row).
Set-TextValue $ws2.Cells.Item(4, 2)  "6048"
Set-TextValue $ws2.Cells.Item(5, 2)  "1565"
Set-TextValue $ws2.Cells.Item(7, 2)  "132"
Set-TextValue $ws2.Cells.Item(8, 2)  "147"
Set-TextValue $ws2.Cells.Item(9, 2)  "202"
Set-TextValue $ws2.Cells.Item(10, 2) "1407"

# New month (Nov) appended at the bottom.
$ws2.Cells.Item(11, 1).Value = 11
Set-TextValue $ws2.Cells.Item(11, 2) "7335"
